$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dollar")
$ws.Activate()

$ws.Range("A3").Value = "P5C"
$ws.Range("A3").Select()
